$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Add the two new DB server references used by the Angular login component
# (local SQL Express instance name, and the localhost alias) below the
# existing "Comandos Consola" table.
$ws.Range("A21").Value = "OFITE-GRUDE8\SQLEXPRESS"
$ws.Range("A22").Value = "localhost"

# Leave the selection where the user's data-entry cursor ended up.
$ws.Range("B22").Select()
